# Outmigration Form feature
# Renames the "In Migration" form to "Outmigration" throughout the workbook,
# tweaks a couple of sheet-view settings, and adds a new boolean flag cell.

$wb = $excel.ActiveWorkbook

$wsColumns  = $wb.Worksheets.Item("columns")
$wsOptions  = $wb.Worksheets.Item("options")
$wsSettings = $wb.Worksheets.Item("settings")

# --- columns sheet ------------------------------------------------------

# 2.1. Migration Type :: pt label
$wsColumns.Cells.Item(5, 9).Value = "2.1. Tipo de Emigração"

# 5.2. Reason for In Migration -> Reason for Outmigration (en/pt)
$wsColumns.Cells.Item(10, 8).Value = "5.2. Reason for Outmigration"
$wsColumns.Cells.Item(10, 9).Value = "5.2. Causas da emigração"

# 5.2.1. Specify other reason for migration :: pt label
$wsColumns.Cells.Item(11, 9).Value = "5.2.1. Especifique outras causas da emigração"

# New boolean cell M4 (mirrors L4/M-column flag already present on neighbouring rows)
$wsColumns.Cells.Item(4, 13).Value = $true

# Sheet view: zoom + selection
$wsColumns.Activate()
$excel.ActiveWindow.Zoom = 190
$wsColumns.Range("A5").Select()

# --- options sheet -------------------------------------------------------

# migtypes row: ENT / Internal InMigration / Imigração Interna -> EXT / External Outmigration / Emigração Externa
$wsOptions.Cells.Item(13, 2).Value = "EXT"
$wsOptions.Cells.Item(13, 3).Value = "External Outmigration"
$wsOptions.Cells.Item(13, 4).Value = "Emigração Externa"

# Sheet view: selection
$wsOptions.Activate()
$wsOptions.Range("B9").Select()

# --- settings sheet -------------------------------------------------------

# form_name::pt now reuses the existing "Emigração Externa" string
$wsSettings.Cells.Item(2, 3).Value = "Emigração Externa"

# Re-activate the columns sheet (it was the originally selected tab)
$wsColumns.Activate()
